$wb = $excel.ActiveWorkbook

$wsGroup = $wb.Worksheets.Item("Group_Members")
$wsChoices = $wb.Worksheets.Item("Choices")

# Add new group members (Group 4)
$wsGroup.Range("A11").Value = 4
$wsGroup.Range("B11").Value = "Austin Nguyen"
$wsGroup.Range("A12").Value = 4
$wsGroup.Range("B12").Value = "Colleen Xu"
$wsGroup.Range("A13").Value = 4
$wsGroup.Range("B13").Value = "Xiao Wang"

# Add new choices for Group 4
$wsChoices.Range("A11").Value = 4
$wsChoices.Range("B11").Value = 1
$wsChoices.Range("C11").Value = 6

$wsChoices.Range("A12").Value = 4
$wsChoices.Range("B12").Value = 2
$wsChoices.Range("C12").Value = 4

$wsChoices.Range("A13").Value = 4
$wsChoices.Range("B13").Value = 3
$wsChoices.Range("C13").Value = 15

# Update selections to reflect new active cells
$wsGroup.Range("B14").Select()
$wsChoices.Range("C14").Select()

# Make Choices the active sheet (tabSelected) and active tab in workbook view
$wsChoices.Activate()

$wb.Save()
